$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3766.348
$ws.Range("J17").Value = 3846.6365
$ws.Range("L17").Value = 11539.9095
$ws.Range("N17").Value = -11875.9095
$ws.Range("H86").Value = 3949.625
$ws.Range("I86").Value = 2349.5
$ws.Range("K86").Value = 2349.5
$ws.Range("M86").Value = -1226.5
$ws.Range("H89").Value = 3949.625
$ws.Range("I89").Value = 2349.5
$ws.Range("K89").Value = 11747.5
$ws.Range("M89").Value = -6131.5
$ws.Range("H113").Value = 3980
$ws.Range("I113").Value = 3980
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3980
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -726
$ws.Range("N113").ClearContents()
$ws.Range("H129").Value = 2999
$ws.Range("J129").Value = 3110
$ws.Range("L129").Value = 9330
$ws.Range("N129").Value = -19330

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3166.3333
$ws.Range("I132").Value = 3249.5
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 9748.5
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -7218.5
$ws.Range("N132").Value = -14060
$ws.Range("H134").Value = 64142.668
$ws.Range("J134").Value = 64142.668
$ws.Range("L134").Value = 64142.668
$ws.Range("N134").Value = -74282.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 242.33333
$ws.Range("I20").Value = 286
$ws.Range("J20").Value = 155
$ws.Range("K20").Value = 286
$ws.Range("L20").Value = 155
$ws.Range("M20").Value = -39
$ws.Range("N20").Value = -649
$ws.Range("H86").Value = 2500
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1377
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2500
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6884
$ws.Range("N89").ClearContents()
$ws.Range("H95").Value = 14031
$ws.Range("J95").Value = 14031
$ws.Range("L95").Value = 14031
$ws.Range("N95").Value = -19523
$ws.Range("H99").Value = 2283.7
$ws.Range("I99").Value = 1976.8572
$ws.Range("K99").Value = 1976.8572
$ws.Range("M99").Value = -478.8571999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 58997.5
$ws.Range("J109").Value = 58997.5
$ws.Range("L109").Value = 58997.5
$ws.Range("N109").Value = -61077.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2000
$ws.Range("J51").Value = 2000
$ws.Range("L51").Value = 6000
$ws.Range("N51").Value = -6920
$ws.Range("H117").Value = 20353.2
$ws.Range("I117").Value = 441.5
$ws.Range("K117").Value = 1324.5
$ws.Range("M117").Value = 2117.5
$ws.Range("H124").Value = 2333
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 2333
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 6999
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -16819
$ws.Range("H129").Value = 1002610.2
$ws.Range("J129").Value = 1670150.4
$ws.Range("L129").Value = 5010451.199999999
$ws.Range("N129").Value = -5020451.199999999
$ws.Range("H130").Value = 1545.1666
$ws.Range("I130").Value = 1465
$ws.Range("J130").Value = 1625.3334
$ws.Range("K130").Value = 4395
$ws.Range("L130").Value = 4876.0002
$ws.Range("M130").Value = 625
$ws.Range("N130").Value = -14916.0002
$ws.Range("H132").Value = 2504.7144
$ws.Range("J132").Value = 2990.8
$ws.Range("L132").Value = 26917.2
$ws.Range("N132").Value = -31977.2
$ws.Range("H139").Value = 2768.75
$ws.Range("I139").Value = 2768.75
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 8306.25
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -3166.25
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8403.237999999999
$ws.Range("I70").Value = 6960.364
$ws.Range("J70").Value = 9990.4
$ws.Range("K70").Value = 6960.364
$ws.Range("L70").Value = 9990.4
$ws.Range("M70").Value = -6690.364
$ws.Range("N70").Value = -10530.4
$ws.Range("H73").Value = 8403.237999999999
$ws.Range("I73").Value = 6960.364
$ws.Range("J73").Value = 9990.4
$ws.Range("K73").Value = 6960.364
$ws.Range("L73").Value = 9990.4
$ws.Range("M73").Value = -6024.364
$ws.Range("N73").Value = -11862.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H16").Value = 3611
$ws.Range("J16").Value = 3611
$ws.Range("L16").Value = 3611
$ws.Range("N16").Value = -3951
$ws.Range("H22").Value = 1629.6364
$ws.Range("I22").Value = 2068.8
$ws.Range("J22").Value = 1263.6666
$ws.Range("K22").Value = 2068.8
$ws.Range("L22").Value = 1263.6666
$ws.Range("M22").Value = -1773.8
$ws.Range("N22").Value = -1853.6666
$ws.Range("H27").Value = 1629.6364
$ws.Range("I27").Value = 2068.8
$ws.Range("J27").Value = 1263.6666
$ws.Range("K27").Value = 2068.8
$ws.Range("L27").Value = 1263.6666
$ws.Range("M27").Value = -1961.8
$ws.Range("N27").Value = -1477.6666
$ws.Range("H61").Value = 2218.158
$ws.Range("I61").Value = 1941.4615
$ws.Range("K61").Value = 1941.4615
$ws.Range("M61").Value = -1739.4615
$ws.Range("H68").Value = 2749.75
$ws.Range("I68").Value = 2499.5
$ws.Range("K68").Value = 2499.5
$ws.Range("M68").Value = -1750.5
$ws.Range("H71").Value = 2749.75
$ws.Range("I71").Value = 2499.5
$ws.Range("K71").Value = 12497.5
$ws.Range("M71").Value = -8753.5
$ws.Range("H98").Value = 19277.5
$ws.Range("J98").Value = 19277.5
$ws.Range("L98").Value = 19277.5
$ws.Range("N98").Value = -25267.5
$ws.Range("H113").Value = 2218.158
$ws.Range("I113").Value = 1941.4615
$ws.Range("K113").Value = 1941.4615
$ws.Range("M113").Value = 228.5385000000001
$ws.Range("H122").Value = 8893.25
$ws.Range("I122").Value = 9710
$ws.Range("J122").Value = 7532
$ws.Range("K122").Value = 29130
$ws.Range("L122").Value = 22596
$ws.Range("M122").Value = -26680
$ws.Range("N122").Value = -27496
$ws.Range("H132").Value = 2584.8667
$ws.Range("I132").Value = 2251.3845
$ws.Range("J132").Value = 4752.5
$ws.Range("K132").Value = 6754.1535
$ws.Range("L132").Value = 14257.5
$ws.Range("M132").Value = -4224.1535
$ws.Range("N132").Value = -19317.5
$ws.Range("H136").Value = 2986.5715
$ws.Range("I136").Value = 2322.4
$ws.Range("J136").Value = 4647
$ws.Range("K136").Value = 6967.200000000001
$ws.Range("L136").Value = 13941
$ws.Range("M136").Value = -4417.200000000001
$ws.Range("N136").Value = -19041

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 57495
$ws.Range("J109").Value = 57495
$ws.Range("L109").Value = 57495
$ws.Range("N109").Value = -60269
$ws.Range("H132").Value = 3142.6296
$ws.Range("J132").Value = 1515.5
$ws.Range("L132").Value = 4546.5
$ws.Range("N132").Value = -9606.5
$ws.Range("H136").Value = 4179.0386
$ws.Range("I136").Value = 4270
$ws.Range("J136").Value = 3974.375
$ws.Range("K136").Value = 12810
$ws.Range("L136").Value = 11923.125
$ws.Range("M136").Value = -10260
$ws.Range("N136").Value = -17023.125
